# Auto commit at 2026-02-25 8:21:59.33
# Update the "Metrics" sheet's monthly values (B2:B13). The "today" sheet
# pulls these via =Metrics!Bn formulas, so it recalculates automatically.

$wb = $excel.ActiveWorkbook
$metrics = $wb.Worksheets.Item("Metrics")

$metrics.Range("B2").Value = 291041.245
$metrics.Range("B3").Value = 266793.69999999995
$metrics.Range("B4").Value = 100815.54000000001
$metrics.Range("B5").Value = 11553
$metrics.Range("B6").Value = 871677.01500000001
$metrics.Range("B7").Value = 719588.15
$metrics.Range("B8").Value = 265250.09999999998
$metrics.Range("B9").Value = 35029
$metrics.Range("B10").Value = 34972928.735000007
$metrics.Range("B11").Value = 32765580.939999998
$metrics.Range("B12").Value = 12211063.960000001
$metrics.Range("B13").Value = 1352936

# Update saved cursor/selection positions to match the authored state.
$metrics.Range("D18").Select() | Out-Null

$today = $wb.Worksheets.Item("today")
$today.Range("F6").Select() | Out-Null
